# Generate Report for Handoff
# Updates the "b.md" rows across the Overview / zh-cn / de-de sheets to
# reflect a fresh handoff (instead of the previous handback) of b.md.

$wb = $excel.ActiveWorkbook

$newZhHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$newDeHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$zhHandoffDate = "2016-08-15 08:52:55"
$deHandoffDate = "2016-08-15 08:53:02"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/25d113ec0d3dfe48c220790b1d3182821ed50085/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/85b064cd59142e2f55bdb31bf6a2bdbeeaf9fa64/e2e/b.md."

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = $deHandoffDate

# --- zh-cn sheet (row 3 = b.md) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
# Write "False" as text (not boolean) via a formula + paste-values round trip.
$zhcn.Range("F3").Formula = '=""&"False"'
$zhcn.Range("F3").Copy()
$zhcn.Range("F3").PasteSpecial(-4163)
$zhcn.Range("G3").Value = $newZhHandoffFile
$zhcn.Range("H3").Value = $zhHandoffDate
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.14

# --- de-de sheet (row 3 = b.md) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Formula = '=""&"False"'
$dede.Range("F3").Copy()
$dede.Range("F3").PasteSpecial(-4163)
$dede.Range("G3").Value = $newDeHandoffFile
$dede.Range("H3").Value = $deHandoffDate
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.14
